# Applies the MINNESOTA_2022 cleanup edit:
#  1. Rename header columns (A1:D1) to snake_case machine-readable names.
#  2. Title-case every Spanish place name in columns A and B (rows 2-981),
#     matching Python's str.title() semantics exactly (incl. the
#     "MonteMorelos" -> "Montemorelos" mid-word-capital fix).
#  3. Nudge the handful of D-column percentage cells that were re-derived
#     with a 1-ulp different floating point result.
#  4. Drop the trailing footnote/source rows (983-987); row 982 is already
#     blank, so the used range collapses back down to A1:D981.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1. Header rename -------------------------------------------------
$ws.Range("A1").Value2 = "mx_state"
$ws.Range("B1").Value2 = "mx_municipality"
$ws.Range("C1").Value2 = "n_matriculas"
$ws.Range("D1").Value2 = "pct_matriculas"

# ---- helper: character is "cased" (has distinct upper/lower form) -----
function IsCasedChar($c) {
    if ($c -eq "") { return $false }
    $up = $c.ToUpper()
    $lo = $c.ToLower()
    if ($up -ne $lo) { return $true }
    return ($c -cmatch '[^\W\d_]')
}

# ---- helper: reproduce Python's str.title() exactly --------------------
function PyTitle($s) {
    $len = $s.Length
    $prevCased = $false
    $result = ""
    for ($i = 0; $i -lt $len; $i++) {
        $ch = $s.Substring($i, 1)
        if (IsCasedChar $ch) {
            if (-not $prevCased) {
                $result = $result + $ch.ToUpper()
            } else {
                $result = $result + $ch.ToLower()
            }
            $prevCased = $true
        } else {
            $result = $result + $ch
            $prevCased = $false
        }
    }
    return $result
}

# ---- 2 & 3. Walk data rows: title-case text, nudge stray float values --
$oldSmall = 0.0009218289085545723
$newSmall = 0.0009218289085545724
$oldBig = 0.009218289085545723
$newBig = 0.009218289085545724

for ($r = 2; $r -le 981; $r++) {
    $aCell = $ws.Cells.Item($r, 1)
    $aVal = $aCell.Value2
    if ($aVal -ne $null -and $aVal -ne "") {
        $aCell.Value2 = PyTitle($aVal)
    }

    $bCell = $ws.Cells.Item($r, 2)
    $bVal = $bCell.Value2
    if ($bVal -ne $null -and $bVal -ne "") {
        $bCell.Value2 = PyTitle($bVal)
    }

    $dCell = $ws.Cells.Item($r, 4)
    $dVal = $dCell.Value2
    if ($dVal -eq $oldSmall) {
        $dCell.Value2 = $newSmall
    } elseif ($dVal -eq $oldBig) {
        $dCell.Value2 = $newBig
    }
}

# ---- 4. Remove the trailing footnote rows 983-987 ----------------------
$ws.Range("A983:A987").EntireRow.Delete()

$wb.Save()
